$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3926408290863037
$ws.Range("E2").Value = 247.7827508240771
$ws.Range("F2").Value = 0.007039267630775069
$ws.Range("G2").Value = 0.006379505252995147
$ws.Range("H2").Value = 0.005992676679445358
$ws.Range("I2").Value = 0.005933050649050004
$ws.Range("J2").Value = 0.005670554378618023
$ws.Range("K2").Value = 0.00554609230832486
$ws.Range("L2").Value = 0.005368068416555367
$ws.Range("M2").Value = 0.005350028676182287
$ws.Range("N2").Value = 0.005330886920923888
$ws.Range("O2").Value = 0.005222996853396169
$ws.Range("P2").Value = 0.005201374725740684
$ws.Range("Q2").Value = 0.005158564432651586
$ws.Range("R2").Value = 0.005158564432651586
$ws.Range("S2").Value = 0.005031695058987714
$ws.Range("T2").Value = 0.005000882620650365
$ws.Range("U2").Value = 0.004997183714750971
$ws.Range("V2").Value = 0.004964961072714322
$ws.Range("W2").Value = 0.004861463602686795
$ws.Range("X2").Value = 0.004856428651974109
$ws.Range("Y2").Value = 0.004830073115479084

$ws.Range("C3").Value = 0.4374735355377197
$ws.Range("E3").Value = 236.3471185846793
$ws.Range("F3").Value = 0.006924576998276322
$ws.Range("G3").Value = 0.006413894287459016
$ws.Range("H3").Value = 0.005695724227595159
$ws.Range("I3").Value = 0.005432667444686222
$ws.Range("J3").Value = 0.00537775210224709
$ws.Range("K3").Value = 0.005198389032611033
$ws.Range("L3").Value = 0.005048080143509221
$ws.Range("M3").Value = 0.004965241438048888
$ws.Range("N3").Value = 0.004965241438048888
$ws.Range("O3").Value = 0.004965241438048888
$ws.Range("P3").Value = 0.004855774544629095
$ws.Range("Q3").Value = 0.004779921238549145
$ws.Range("R3").Value = 0.004730227818857203
$ws.Range("S3").Value = 0.004730227818857203
$ws.Range("T3").Value = 0.004707253353470391
$ws.Range("U3").Value = 0.00467207040456914
$ws.Range("V3").Value = 0.00466485475020488
$ws.Range("W3").Value = 0.004639791634209141
$ws.Range("X3").Value = 0.004631921710464042
$ws.Range("Y3").Value = 0.004607156307693551

$ws.Range("C4").Value = 0.5156247615814209
$ws.Range("E4").Value = 235.6774715577303
$ws.Range("F4").Value = 0.006991883289192024
$ws.Range("G4").Value = 0.006216056926719194
$ws.Range("H4").Value = 0.006100775546888117
$ws.Range("I4").Value = 0.005716206943185562
$ws.Range("J4").Value = 0.005643899836768592
$ws.Range("K4").Value = 0.005576468286968062
$ws.Range("L4").Value = 0.005389589799091743
$ws.Range("M4").Value = 0.005251061208827164
$ws.Range("N4").Value = 0.005232184215597051
$ws.Range("O4").Value = 0.005062461164214159
$ws.Range("P4").Value = 0.005039260274946298
$ws.Range("Q4").Value = 0.004948983830580043
$ws.Range("R4").Value = 0.004813853858165856
$ws.Range("S4").Value = 0.00471486728055802
$ws.Range("T4").Value = 0.00471486728055802
$ws.Range("U4").Value = 0.00471486728055802
$ws.Range("V4").Value = 0.004680508199265471
$ws.Range("W4").Value = 0.004626800425830808
$ws.Range("X4").Value = 0.004612829133752524
$ws.Range("Y4").Value = 0.004594102759409946

$ws.Range("C5").Value = 0.4218757152557373
$ws.Range("E5").Value = 234.8336414771766
$ws.Range("F5").Value = 0.007164587072448732
$ws.Range("G5").Value = 0.006050059988423096
$ws.Range("H5").Value = 0.005756281780270556
$ws.Range("I5").Value = 0.005593430202232652
$ws.Range("J5").Value = 0.005450650160137281
$ws.Range("K5").Value = 0.005439892877537763
$ws.Range("L5").Value = 0.005262587291965731
$ws.Range("M5").Value = 0.005237266037350487
$ws.Range("N5").Value = 0.004952001170737065
$ws.Range("O5").Value = 0.004917439462269275
$ws.Range("P5").Value = 0.004909506667708198
$ws.Range("Q5").Value = 0.004810135280202383
$ws.Range("R5").Value = 0.004740524163207559
$ws.Range("S5").Value = 0.004697780683187697
$ws.Range("T5").Value = 0.004682468173245597
$ws.Range("U5").Value = 0.004655913975576752
$ws.Range("V5").Value = 0.004642109381614089
$ws.Range("W5").Value = 0.00459128954086341
$ws.Range("X5").Value = 0.00459128954086341
$ws.Range("Y5").Value = 0.004577653829964454

$ws.Range("C6").Value = 0.4375348091125488
$ws.Range("E6").Value = 241.3709982473665
$ws.Range("F6").Value = 0.006903139592921585
$ws.Range("G6").Value = 0.006204056699690823
$ws.Range("H6").Value = 0.005757271650882048
$ws.Range("I6").Value = 0.005698095064417109
$ws.Range("J6").Value = 0.005554323254832533
$ws.Range("K6").Value = 0.005554323254832533
$ws.Range("L6").Value = 0.005546628031019561
$ws.Range("M6").Value = 0.005391319488743227
$ws.Range("N6").Value = 0.005391319488743227
$ws.Range("O6").Value = 0.005244614828560066
$ws.Range("P6").Value = 0.00515059944027896
$ws.Range("Q6").Value = 0.004947169500467859
$ws.Range("R6").Value = 0.004947169500467859
$ws.Range("S6").Value = 0.004892568588757074
$ws.Range("T6").Value = 0.004858157973933691
$ws.Range("U6").Value = 0.004825713092965233
$ws.Range("V6").Value = 0.004785250813344502
$ws.Range("W6").Value = 0.004785250813344502
$ws.Range("X6").Value = 0.004719877406016137
$ws.Range("Y6").Value = 0.004705087685133849

$ws.Range("C7").Value = 0.5569784641265869
$ws.Range("E7").Value = 242.6358888408104
$ws.Range("F7").Value = 0.006887718867366348
$ws.Range("G7").Value = 0.006063487179328527
$ws.Range("H7").Value = 0.005877852989776113
$ws.Range("I7").Value = 0.005603297877380177
$ws.Range("J7").Value = 0.005560686992261497
$ws.Range("K7").Value = 0.005359086019515882
$ws.Range("L7").Value = 0.005233873793228222
$ws.Range("M7").Value = 0.005146086463740425
$ws.Range("N7").Value = 0.005146086463740425
$ws.Range("O7").Value = 0.005142828519656424
$ws.Range("P7").Value = 0.004995613644830909
$ws.Range("Q7").Value = 0.004937869913442625
$ws.Range("R7").Value = 0.004870360548434955
$ws.Range("S7").Value = 0.004870360548434955
$ws.Range("T7").Value = 0.004854372168062969
$ws.Range("U7").Value = 0.004786134759725297
$ws.Range("V7").Value = 0.004757437572317716
$ws.Range("W7").Value = 0.004757437572317716
$ws.Range("X7").Value = 0.004753879162830339
$ws.Range("Y7").Value = 0.004729744421848154

$ws.Range("C8").Value = 0.4241547584533691
$ws.Range("E8").Value = 251.4473630892771
$ws.Range("F8").Value = 0.007109132424406193
$ws.Range("G8").Value = 0.006042823208021887
$ws.Range("H8").Value = 0.00592215226963446
$ws.Range("I8").Value = 0.005632143238663029
$ws.Range("J8").Value = 0.005632143238663029
$ws.Range("K8").Value = 0.005520109195946438
$ws.Range("L8").Value = 0.005520109195946438
$ws.Range("M8").Value = 0.005520109195946438
$ws.Range("N8").Value = 0.005520109195946438
$ws.Range("O8").Value = 0.005296690572315028
$ws.Range("P8").Value = 0.005296690572315028
$ws.Range("Q8").Value = 0.005205694379888548
$ws.Range("R8").Value = 0.00513515907119778
$ws.Range("S8").Value = 0.005065338725895311
$ws.Range("T8").Value = 0.004969108437109032
$ws.Range("U8").Value = 0.004969108437109032
$ws.Range("V8").Value = 0.004948391964720776
$ws.Range("W8").Value = 0.004947991265724102
$ws.Range("X8").Value = 0.004901508052422556
$ws.Range("Y8").Value = 0.004901508052422556

$ws.Range("C9").Value = 0.4375081062316895
$ws.Range("E9").Value = 238.6743466546432
$ws.Range("F9").Value = 0.007110440230760829
$ws.Range("G9").Value = 0.006290415420899859
$ws.Range("H9").Value = 0.005660798712296734
$ws.Range("I9").Value = 0.005660798712296734
$ws.Range("J9").Value = 0.005526496296618408
$ws.Range("K9").Value = 0.005278738572034594
$ws.Range("L9").Value = 0.005238727854869574
$ws.Range("M9").Value = 0.005080838014660986
$ws.Range("N9").Value = 0.005080838014660986
$ws.Range("O9").Value = 0.005010461058546392
$ws.Range("P9").Value = 0.004962920583519496
$ws.Range("Q9").Value = 0.004887595201686003
$ws.Range("R9").Value = 0.004887595201686003
$ws.Range("S9").Value = 0.004887595201686003
$ws.Range("T9").Value = 0.00483145318522345
$ws.Range("U9").Value = 0.004766390292302868
$ws.Range("V9").Value = 0.004715112274485579
$ws.Range("W9").Value = 0.004692808336503287
$ws.Range("X9").Value = 0.004664200919009476
$ws.Range("Y9").Value = 0.004652521377283492

$ws.Range("C10").Value = 0.4531004428863525
$ws.Range("E10").Value = 241.070044650065
$ws.Range("F10").Value = 0.006982259873881369
$ws.Range("G10").Value = 0.006023874135017891
$ws.Range("H10").Value = 0.006018370504520538
$ws.Range("I10").Value = 0.005713122579205059
$ws.Range("J10").Value = 0.005708302549488474
$ws.Range("K10").Value = 0.005377935670662656
$ws.Range("L10").Value = 0.005377935670662656
$ws.Range("M10").Value = 0.005083272663345992
$ws.Range("N10").Value = 0.005083272663345992
$ws.Range("O10").Value = 0.005076052182641191
$ws.Range("P10").Value = 0.005002875746391601
$ws.Range("Q10").Value = 0.004910801382794216
$ws.Range("R10").Value = 0.004897457529354403
$ws.Range("S10").Value = 0.004817160699171188
$ws.Range("T10").Value = 0.004817160699171188
$ws.Range("U10").Value = 0.004798742237798486
$ws.Range("V10").Value = 0.004768103586381387
$ws.Range("W10").Value = 0.004745466125138997
$ws.Range("X10").Value = 0.004701409642272272
$ws.Range("Y10").Value = 0.00469922114327612

$ws.Range("C11").Value = 0.4219002723693848
$ws.Range("E11").Value = 244.0361624010875
$ws.Range("F11").Value = 0.006939853198287947
$ws.Range("G11").Value = 0.006184302063529063
$ws.Range("H11").Value = 0.005863156135026648
$ws.Range("I11").Value = 0.005706046768763199
$ws.Range("J11").Value = 0.005393229082475375
$ws.Range("K11").Value = 0.005352543038583742
$ws.Range("L11").Value = 0.005352543038583742
$ws.Range("M11").Value = 0.00502433553759918
$ws.Range("N11").Value = 0.00502433553759918
$ws.Range("O11").Value = 0.00502433553759918
$ws.Range("P11").Value = 0.004993050516996236
$ws.Range("Q11").Value = 0.004973710896965481
$ws.Range("R11").Value = 0.004973710896965481
$ws.Range("S11").Value = 0.004852222514633417
$ws.Range("T11").Value = 0.004852222514633417
$ws.Range("U11").Value = 0.004852222514633417
$ws.Range("V11").Value = 0.004781483924785507
$ws.Range("W11").Value = 0.004781483924785507
$ws.Range("X11").Value = 0.004768052431868303
$ws.Range("Y11").Value = 0.004757040202750243
